$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 11760.071
$ws.Range("I80").Value = 17750.5
$ws.Range("J80").Value = 7267.25
$ws.Range("K80").Value = 53251.5
$ws.Range("L80").Value = 21801.75
$ws.Range("M80").Value = -52253.5
$ws.Range("N80").Value = -23797.75

$ws.Range("H83").Value = 11760.071
$ws.Range("I83").Value = 17750.5
$ws.Range("J83").Value = 7267.25
$ws.Range("K83").Value = 159754.5
$ws.Range("L83").Value = 65405.25
$ws.Range("M83").Value = -154762.5
$ws.Range("N83").Value = -75389.25

$ws.Range("H137").Value = 2881.476
$ws.Range("I137").Value = 2970.6072
$ws.Range("J137").Value = 2703.2144
$ws.Range("K137").Value = 8911.821599999999
$ws.Range("L137").Value = 8109.6432
$ws.Range("M137").Value = -6361.821599999999
$ws.Range("N137").Value = -13209.6432

$ws.Range("H138").Value = 6850.9272
$ws.Range("I138").Value = 4238.8667
$ws.Range("J138").Value = 7334.642
$ws.Range("K138").Value = 12716.6001
$ws.Range("L138").Value = 22003.926
$ws.Range("M138").Value = -7576.6001
$ws.Range("N138").Value = -32283.926

$ws.Range("H140").Value = 55356.215
$ws.Range("J140").Value = 55356.215
$ws.Range("L140").Value = 55356.215
$ws.Range("N140").Value = -65716.215

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 549.5
$ws.Range("I4").Value = 152.5
$ws.Range("J4").Value = 748
$ws.Range("K4").Value = 152.5
$ws.Range("L4").Value = 748
$ws.Range("M4").Value = -36.5
$ws.Range("N4").Value = -980

$ws.Range("H28").Value = 8243.75
$ws.Range("J28").Value = 4975
$ws.Range("L28").Value = 4975
$ws.Range("N28").Value = -5359

$ws.Range("H32").Value = 8544.84
$ws.Range("I32").Value = 4863.4414
$ws.Range("J32").Value = 20869.521
$ws.Range("K32").Value = 4863.4414
$ws.Range("L32").Value = 20869.521
$ws.Range("M32").Value = -4576.4414
$ws.Range("N32").Value = -21443.521

$ws.Range("H70").Value = 50000
$ws.Range("J70").Value = 50000
$ws.Range("L70").Value = 50000
$ws.Range("N70").Value = -50540

$ws.Range("H73").Value = 50000
$ws.Range("J73").Value = 50000
$ws.Range("L73").Value = 50000
$ws.Range("N73").Value = -51872

$ws.Range("H99").Value = 8243.75
$ws.Range("J99").Value = 4975
$ws.Range("L99").Value = 4975
$ws.Range("N99").Value = -10965

$ws.Range("H132").Value = 5921.4473
$ws.Range("I132").Value = 5193.5806
$ws.Range("J132").Value = 9144.857
$ws.Range("K132").Value = 15580.7418
$ws.Range("L132").Value = 27434.571
$ws.Range("M132").Value = -13050.7418
$ws.Range("N132").Value = -32494.571

$ws.Range("H135").Value = 92000
$ws.Range("J135").Value = 92000
$ws.Range("L135").Value = 92000
$ws.Range("N135").Value = -102140

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1922.9474
$ws.Range("I20").Value = 2083.9167
$ws.Range("J20").Value = 1647
$ws.Range("K20").Value = 2083.9167
$ws.Range("L20").Value = 1647
$ws.Range("M20").Value = -1836.9167
$ws.Range("N20").Value = -2141

$ws.Range("H28").Value = 21250
$ws.Range("J28").Value = 21250
$ws.Range("L28").Value = 21250
$ws.Range("N28").Value = -21838

$ws.Range("H105").Value = 94108.73
$ws.Range("I105").Value = 103119.7
$ws.Range("K105").Value = 103119.7
$ws.Range("M105").Value = -101372.7

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("H118").Value = 70000
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 70000
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 70000
$ws.Range("M118").ClearContents()
$ws.Range("N118").Value = -73314

$ws.Range("H134").Value = 24566.334
$ws.Range("I134").Value = 3812.4255
$ws.Range("K134").Value = 11437.2765
$ws.Range("M134").Value = -8902.2765

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 546.2308
$ws.Range("I7").Value = 550.1818
$ws.Range("J7").Value = 524.5
$ws.Range("K7").Value = 550.1818
$ws.Range("L7").Value = 524.5
$ws.Range("M7").Value = -437.1818
$ws.Range("N7").Value = -750.5

$ws.Range("H31").Value = 67377.06
$ws.Range("I31").Value = 1890
$ws.Range("J31").Value = 141050
$ws.Range("K31").Value = 1890
$ws.Range("L31").Value = 141050
$ws.Range("M31").Value = -1595
$ws.Range("N31").Value = -141640

$ws.Range("H34").Value = 67377.06
$ws.Range("I34").Value = 1890
$ws.Range("J34").Value = 141050
$ws.Range("K34").Value = 1890
$ws.Range("L34").Value = 141050
$ws.Range("M34").Value = -1688
$ws.Range("N34").Value = -141454

$ws.Range("H68").Value = 34000
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 34000
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H132").Value = 1771.641
$ws.Range("I132").Value = 1301.3
$ws.Range("J132").Value = 3339.4443
$ws.Range("K132").Value = 3903.9
$ws.Range("L132").Value = 10018.3329
$ws.Range("M132").Value = -1373.9
$ws.Range("N132").Value = -15078.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1611570
$ws.Range("I4").Value = 1952479.8
$ws.Range("K4").Value = 5857439.4
$ws.Range("M4").Value = -5857327.4

$ws.Range("H18").Value = 953.5714
$ws.Range("I18").Value = 653.73334
$ws.Range("K18").Value = 1961.20002
$ws.Range("M18").Value = -1792.20002

$ws.Range("H92").Value = 1401.1765
$ws.Range("I92").Value = 1199.8
$ws.Range("J92").Value = 1485.0834
$ws.Range("K92").Value = 3599.4
$ws.Range("L92").Value = 4455.2502
$ws.Range("M92").Value = -2351.4
$ws.Range("N92").Value = -6951.2502

$ws.Range("H131").Value = 65985.31
$ws.Range("I131").Value = 68508.39999999999
$ws.Range("J131").Value = 63759.06
$ws.Range("K131").Value = 205525.2
$ws.Range("L131").Value = 191277.18
$ws.Range("M131").Value = -200485.2
$ws.Range("N131").Value = -201357.18

$ws.Range("H139").Value = 4266.4614
$ws.Range("I139").Value = 3229.7778
$ws.Range("J139").Value = 6599
$ws.Range("K139").Value = 9689.3334
$ws.Range("L139").Value = 19797
$ws.Range("M139").Value = -4549.3334
$ws.Range("N139").Value = -30077

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 33834.6
$ws.Range("J101").Value = 33834.6
$ws.Range("L101").Value = 33834.6
$ws.Range("N101").Value = -40324.6

$ws.Range("H132").Value = 42639.574
$ws.Range("I132").Value = 12528
$ws.Range("J132").Value = 102862.73
$ws.Range("K132").Value = 37584
$ws.Range("L132").Value = 308588.19
$ws.Range("M132").Value = -35054
$ws.Range("N132").Value = -313648.19

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 266.69232
$ws.Range("I16").Value = 253.36
$ws.Range("J16").Value = 600
$ws.Range("K16").Value = 253.36
$ws.Range("L16").Value = 600
$ws.Range("M16").Value = -83.36000000000001
$ws.Range("N16").Value = -940

$ws.Range("H46").Value = 3041.8333
$ws.Range("I46").Value = 2855.5557
$ws.Range("K46").Value = 2855.5557
$ws.Range("M46").Value = -2667.5557

$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()

$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("M124").ClearContents()

$ws.Range("H132").Value = 4328.75
$ws.Range("I132").Value = 4217
$ws.Range("J132").Value = 5111
$ws.Range("K132").Value = 12651
$ws.Range("L132").Value = 15333
$ws.Range("M132").Value = -10121
$ws.Range("N132").Value = -20393

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 60092.824
$ws.Range("I96").Value = 72652
$ws.Range("J96").Value = 1483.3334
$ws.Range("K96").Value = 72652
$ws.Range("L96").Value = 1483.3334
$ws.Range("M96").Value = -71279
$ws.Range("N96").Value = -4229.3334

$ws.Range("H136").Value = 330170.1
$ws.Range("I136").Value = 301247.22
$ws.Range("J136").Value = 667603.3
$ws.Range("K136").Value = 903741.6599999999
$ws.Range("L136").Value = 2002809.9
$ws.Range("M136").Value = -901191.6599999999
$ws.Range("N136").Value = -2007909.9
